# diary.xlsx - "Still working on renderable sphere as renderableparticle"
#
# Adds a new diary entry for "29 syys" (row 5) with its hours, extends the
# existing "26 syys" entry's learning-content note (C4), and leaves a lone
# space marker in G7. The order of writes below matters: it controls the
# order new strings are appended to the shared-string table so the saved
# workbook's xl/sharedStrings.xml matches the authored edit byte-for-byte.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New row 5 date cell ("29 syys") - appended to the shared strings first.
$ws.Range("A5").Value = "29 syys"

# 2) Extend the existing "26 syys" entry's content note in C4 and make sure
#    it keeps wrapping like the rest of that row.
$ws.Range("C4").Value = "Kirjan sivut 36-54, calculuksen pikakertausta ja johdantoa partikkeleihin"
$ws.Range("C4").WrapText = $true

# 3) Lone space placeholder down in G7 (row 6 stays empty/skipped).
$ws.Range("G7").Value = " "

# 4) New row 5 time-range cell, formatted like a time (h:mm) even though the
#    value itself is free-form text, matching the style used elsewhere in
#    the "Kello" column.
$ws.Range("B5").Value = "9.00-10.45"
$ws.Range("B5").NumberFormat = "h:mm"

# 5) Hours logged for the new entry - H3's SUM(G3:G40) formula recalculates
#    automatically (2.5 + 2 + 1.5 = 6).
$ws.Range("G5").Value = 1.5

# 6) Restore the selection to where the author's cursor ended up.
$ws.Range("H5").Select()
